# Update countries & provincias Spain
# - Refresh the "Datos actualizados ..." timestamp in A1
# - Update several countries' stats to the latest figures
# - Because the sheet is kept sorted by "Casos totales" (col B) descending,
#   a couple of countries leapfrog their neighbours now that their totals
#   changed (Singapur overtakes Bielorrusia; Moldavia overtakes Croacia and
#   Marruecos), so those rows are rewritten in their new order/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 16:52"

# Estados Unidos (row 4)
Set-Row 4 @("Estados Unidos", 505237, 2361, 28088, 458299, 10947, 103, 18850)

# Suiza (row 14)
Set-Row 14 @("Suiza", 24900, 349, 11100, 12785, 386, 13, 1015)

# Brasil (row 17)
Set-Row 17 @("Brasil", 19967, 178, 173, 18719, 296, 7, 1075)

# Singapur / Bielorrusia swap places (rows 52-53) with updated figures
Set-Row 52 @("Singapur", 2299, 191, 528, 1763, 31, 1, 8)
Set-Row 53 @("Bielorrusia", 2226, 245, 172, 2031, 72, 4, 23)

# Moldavia moves ahead of Croacia and Marruecos (rows 60-62) with updated figures
Set-Row 60 @("Moldavia", 1560, 122, 75, 1455, 80, 1, 30)
Set-Row 61 @("Croacia", 1534, 39, 323, 1190, 32, 0, 21)
Set-Row 62 @("Marruecos", 1527, 79, 141, 1276, 1, 3, 110)

# Armenia (row 73)
Set-Row 73 @("Armenia", 977, 40, 173, 791, 30, 1, 13)

# Bulgaria (row 82)
Set-Row 82 @("Bulgaria", 661, 26, 62, 571, 32, 3, 28)

# Mauricio (row 105)
Set-Row 105 @("Mauricio", 319, 1, 28, 282, 3, 0, 9)
